$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'43.015.17"

$ws.Range("D3").Formula = "'2.237.13"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Formula = "'246.30"
$ws.Range("E5").Value = "  +3.84%  "

$ws.Range("D6").Formula = "'0.619"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").Formula = "'75.57"
$ws.Range("E7").Value = "  +7.36%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Formula = "'0.612"
$ws.Range("E9").Value = "  +5.25%  "

$ws.Range("D10").Formula = "'40.84"
$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("D11").Formula = "'0.0932"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").Formula = "'55.48"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").Formula = "'6.94"
$ws.Range("E13").Value = "  +2.25%  "

$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D15").Formula = "'2.556.66"
$ws.Range("E15").Value = "  +2.17%  "

$ws.Range("D16").Formula = "'14.67"
$ws.Range("E16").Value = "  +5.14%  "

$ws.Range("D17").Formula = "'2.237.92"
$ws.Range("E17").Value = "  +2.37%  "

$ws.Range("D18").Formula = "'0.812"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Formula = "'42.944.90"
$ws.Range("E19").Value = "  +4.57%  "

$ws.Range("D20").Formula = "'0.0000104"
$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("D21").Formula = "'71.04"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").Formula = "'5.99"
$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").Formula = "'10.37"
$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Formula = "'230.42"
$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Formula = "'2.19"
$ws.Range("E25").Value = "  +11.62%  "

$ws.Range("D27").Formula = "'10.92"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("E28").Value = "  -5.35%  "

$ws.Range("E29").Value = "  +1.52%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Formula = "'2.21"
$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("D31").Formula = "'173.90"
$ws.Range("E31").Value = "  +3.93%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Formula = "'37.26"
$ws.Range("E32").Value = "  +19.68%  "

$ws.Range("D33").Formula = "'20.31"
$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("E34").Value = "  +2.92%  "

$ws.Range("D35").Formula = "'5.35"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("D36").Formula = "'0.122"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  +7.38%  "

$ws.Range("D38").Formula = "'4.35"
$ws.Range("E38").Value = "  +5.35%  "

$ws.Range("E39").Value = "  +16.47%  "

$ws.Range("D40").Formula = "'13.04"
$ws.Range("E40").Value = "  +6.21%  "

$ws.Range("E41").Value = "  +2.59%  "

$ws.Range("D42").Formula = "'5.58"
$ws.Range("E42").Value = "  +2.61%  "

$ws.Range("D43").Formula = "'0.199"
$ws.Range("E43").Value = "  +4.20%  "

$ws.Range("D44").Formula = "'60.02"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("D45").Formula = "'105.40"
$ws.Range("E45").Value = "  +7.01%  "

$ws.Range("D46").Formula = "'8.61"
$ws.Range("E46").Value = "  +3.15%  "

$ws.Range("D47").Formula = "'0.0990"
$ws.Range("E47").Value = "  +1.61%  "

$ws.Range("D48").Formula = "'0.442"
$ws.Range("E48").Value = "  +19.92%  "

$ws.Range("D49").Formula = "'1.10"
$ws.Range("E49").Value = "  +1.18%  "

$ws.Range("D50").Formula = "'2.30"
$ws.Range("E50").Value = "  +2.53%  "

$ws.Range("E51").Value = "  +0.97%  "
